$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet from "Hoja1" to "General"
$ws.Name = "General"

# Move the active selection from C31 to B30
$ws.Range("B30").Select()
